# Update "想去人数" (interest count) values in the F column of the
# 展览 (sheet1) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$sheet1Updates = @{
    2  = 7294
    3  = 23
    5  = 27
    6  = 571
    7  = 198
    8  = 139
    12 = 228
    14 = 469
    16 = 1877
    19 = 3821
    21 = 253
    23 = 47
    26 = 2508
    28 = 327
    30 = 8
    31 = 47
    32 = 10
    33 = 29
    38 = 39
    39 = 1497
    40 = 165
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型": row -> new F value
$sheet4Updates = @{
    2  = 7294
    3  = 23
    5  = 27
    7  = 571
    8  = 198
    9  = 139
    13 = 228
    15 = 469
    17 = 1877
    20 = 3821
    22 = 253
    24 = 47
    27 = 2508
    29 = 327
    31 = 8
    32 = 47
    33 = 10
    34 = 29
    39 = 39
    40 = 1497
    41 = 165
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
